$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.747.27"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.875.08"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.27"
$ws.Range("E5").Value = "  +6.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.51"
$ws.Range("E6").Value = "  -3.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  -2.33%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.715"
$ws.Range("E9").Value = "  -2.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  -5.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000324"
$ws.Range("E11").Value = "  -7.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.82"
$ws.Range("E12").Value = "  -3.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.44"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.485.64"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.46"
$ws.Range("E15").Value = "  +7.88%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.16"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.808.93"
$ws.Range("E17").Value = "  -3.63%  "

$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.20"
$ws.Range("E19").Value = "  +2.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.695.57"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "418.06"
$ws.Range("E21").Value = "  -5.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("E22").Value = "  +1.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.09"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.10"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("E25").Value = "  +7.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -2.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.50"
$ws.Range("E27").Value = "  -5.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.62"
$ws.Range("E28").Value = "  -4.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.50"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "680.05"
$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.93"
$ws.Range("E31").Value = "  +14.22%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").Value = "  -4.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.85"
$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "67.06"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.446"
$ws.Range("E35").Value = "  -5.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0849"
$ws.Range("E36").Value = "  -8.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.70"
$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +14.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0476"
$ws.Range("E42").Value = "  -2.91%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +4.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.85"
$ws.Range("E44").Value = "  -3.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("E46").Value = "  -1.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("E48").Value = "  +13.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.29"
$ws.Range("E49").Value = "  -2.83%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.82"
$ws.Range("E50").Value = "  +3.89%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.90"
$ws.Range("E51").Value = "  -0.60%  "

